$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("CompStat_1")

# --- Header text updates ---
# These shared strings are rich-text runs split across several <r> elements;
# writing the full combined text collapses them into a single plain run (same visible text/style).
$ws.Range("A8").Value = "Volume 30   Number  50"
$ws.Range("C9").Value = "Report Covering the Week  12/11/2023  Through  12/17/2023"

# --- Crime statistics table updates (rows 14-30) ---
# For cells whose value flips from numeric to the text placeholders "0"/"***.*",
# copy from an existing untouched donor cell (row 14) so value AND style (s=14) match exactly.
$ws.Range("N14").Value = -91.304347826087
$ws.Range("C14").Copy($ws.Range("C15"))
$ws.Range("C14").Copy($ws.Range("G15"))
$ws.Range("E14").Copy($ws.Range("H15"))
$ws.Range("C16").Value = 4
$ws.Range("C14").Copy($ws.Range("D16"))
$ws.Range("E14").Copy($ws.Range("E16"))
$ws.Range("F16").Value = 17
$ws.Range("G16").Value = 5
$ws.Range("H16").Value = 240
$ws.Range("I16").Value = 133
$ws.Range("K16").Value = -10.135135135135
$ws.Range("L16").Value = -5
$ws.Range("M16").Value = -28.108108108108
$ws.Range("N16").Value = -84.552845528455
$ws.Range("C17").Value = 2
$ws.Range("D17").Value = 3
$ws.Range("E17").Value = -33.333333333333
$ws.Range("F17").Value = 13
$ws.Range("G17").Value = 12
$ws.Range("H17").Value = 8.333333333333
$ws.Range("I17").Value = 174
$ws.Range("J17").Value = 180
$ws.Range("K17").Value = -3.333333333333
$ws.Range("L17").Value = -4.395604395604
$ws.Range("M17").Value = 74
$ws.Range("N17").Value = -61.674008810572
$ws.Range("C18").Value = 2
$ws.Range("D18").Value = 4
$ws.Range("E18").Value = -50
$ws.Range("G18").Value = 13
$ws.Range("H18").Value = -30.76923076923
$ws.Range("I18").Value = 165
$ws.Range("J18").Value = 177
$ws.Range("K18").Value = -6.779661016949
$ws.Range("L18").Value = 26.923076923076
$ws.Range("M18").Value = 32
$ws.Range("N18").Value = -84.119345524542
$ws.Range("C19").Value = 10
$ws.Range("D19").Value = 7
$ws.Range("E19").Value = 42.857142857142
$ws.Range("F19").Value = 29
$ws.Range("G19").Value = 39
$ws.Range("H19").Value = -25.641025641025
$ws.Range("I19").Value = 494
$ws.Range("J19").Value = 512
$ws.Range("K19").Value = -3.515625
$ws.Range("L19").Value = 7.158351409978
$ws.Range("M19").Value = -4.816955684007
$ws.Range("N19").Value = -54.343807763401
$ws.Range("C14").Copy($ws.Range("C20"))
$ws.Range("D20").Value = 4
$ws.Range("E20").Value = -100
$ws.Range("F20").Value = 5
$ws.Range("G20").Value = 7
$ws.Range("H20").Value = -28.571428571428
$ws.Range("J20").Value = 86
$ws.Range("K20").Value = 10.465116279069
$ws.Range("L20").Value = 11.764705882352
$ws.Range("M20").Value = 120.93023255814
$ws.Range("N20").Value = -90.165631469979
$ws.Range("C21").Value = 18
$ws.Range("D21").Value = 18
$ws.Range("E21").Value = 0
$ws.Range("F21").Value = 75
$ws.Range("G21").Value = 76
$ws.Range("H21").Value = -1.315789473684
$ws.Range("I21").Value = 1075
$ws.Range("J21").Value = 1113
$ws.Range("K21").Value = -3.414195867026
$ws.Range("L21").Value = 6.120434353405
$ws.Range("M21").Value = 9.137055837563
$ws.Range("N21").Value = -75.896860986547
$ws.Range("C14").Copy($ws.Range("D22"))
$ws.Range("E14").Copy($ws.Range("E22"))
$ws.Range("G22").Value = 3
$ws.Range("H22").Value = 33.333333333333
$ws.Range("L22").Value = -7.407407407407
$ws.Range("C23").Value = 2
$ws.Range("E23").Value = 100
$ws.Range("F23").Value = 7
$ws.Range("G23").Value = 7
$ws.Range("H23").Value = 0
$ws.Range("J23").Value = 120
$ws.Range("K23").Value = -10.833333333333
$ws.Range("L23").Value = 0.943396226415
$ws.Range("C24").Value = 30
$ws.Range("D24").Value = 26
$ws.Range("E24").Value = 15.384615384615
$ws.Range("F24").Value = 89
$ws.Range("G24").Value = 116
$ws.Range("H24").Value = -23.275862068965
$ws.Range("I24").Value = 1502
$ws.Range("J24").Value = 1789
$ws.Range("K24").Value = -16.042481833426
$ws.Range("L24").Value = 19.39586645469
$ws.Range("M24").Value = 45.542635658914
$ws.Range("C25").Value = 11
$ws.Range("D25").Value = 8
$ws.Range("E25").Value = 37.5
$ws.Range("F25").Value = 24
$ws.Range("H25").Value = 26.315789473684
$ws.Range("I25").Value = 297
$ws.Range("J25").Value = 296
$ws.Range("K25").Value = 0.337837837837
$ws.Range("L25").Value = 10.408921933085
$ws.Range("M25").Value = -6.896551724137
$ws.Range("C14").Copy($ws.Range("C26"))
$ws.Range("C14").Copy($ws.Range("G26"))
$ws.Range("E14").Copy($ws.Range("H26"))
$ws.Range("F27").Value = 4
$ws.Range("G27").Value = 2
$ws.Range("H27").Value = 100
$ws.Range("I27").Value = 43
$ws.Range("K27").Value = -15.686274509803
$ws.Range("L27").Value = -17.307692307692
$ws.Range("L28").Value = -30
$ws.Range("N28").Value = -86.792452830188
$ws.Range("L29").Value = 0
$ws.Range("N29").Value = -88
$ws.Range("F30").Value = 2
